# Calendar feature updated with teacher-layout load
# Append additional scenario rows (31-44) describing the new
# "Planning a course" / "Monitorizing weekly and monthly basis courses
# in Calendar Feature" teacher-layout test results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Planning a course", "PASSED", "chrome"),
    @("Planning a course", "PASSED", "chrome"),
    @("Monitorizing weekly and monthly basis courses in Calendar Feature", "FAILED", "chrome"),
    @("Planning a course", "FAILED", "chrome"),
    @("Monitorizing weekly and monthly basis courses in Calendar Feature", "FAILED", "chrome"),
    @("Planning a course", "PASSED", "chrome"),
    @("Monitorizing weekly and monthly basis courses in Calendar Feature", "PASSED", "chrome"),
    @("Monitorizing weekly and monthly basis courses in Calendar Feature", "PASSED", "chrome"),
    @("Planning a course", "FAILED", "chrome"),
    @("Planning a course", "FAILED", "chrome"),
    @("Planning a course", "FAILED", "chrome"),
    @("Monitorizing weekly and monthly basis courses in Calendar Feature", "FAILED", "chrome"),
    @("Planning a course", "PASSED", "chrome"),
    @("Monitorizing weekly and monthly basis courses in Calendar Feature", "PASSED", "chrome")
)

$startRow = 31
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
